$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview" (sheet1) ----
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "ffff5977e3da-bc1e-43b4-b9a8-4e91e21a3592.md"
$ws1.Range("A3").Value = "ffffff9fb0db73-d95f-4e7d-9243-98fa0070a285.md"
$ws1.Range("A4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

# ---- Sheet "zh-cn" (sheet2) ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "ffff5977e3da-bc1e-43b4-b9a8-4e91e21a3592.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-02-17 03:21:44"
$ws2.Range("E2").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.md"
$ws2.Range("F2").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-02-17 03:22:25"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "ffffff9fb0db73-d95f-4e7d-9243-98fa0070a285.md"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-02-17 03:21:44"
$ws2.Range("E3").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.md"
$ws2.Range("F3").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-02-17 03:22:25"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.34df0f84efa86276470169becf660a08a05a939e.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-02-17 03:25:44"
$ws2.Range("E4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.md"
$ws2.Range("F4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.34df0f84efa86276470169becf660a08a05a939e.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-02-17 03:24:43"
$ws2.Range("H4").Value = "Include"

# ---- Sheet "de-de" (sheet3) ----
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "ffff5977e3da-bc1e-43b4-b9a8-4e91e21a3592.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.de-de.xlf"
$ws3.Range("D2").Value = "2016-02-17 03:21:54"
$ws3.Range("E2").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.md"
$ws3.Range("F2").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.de-de.xlf"
$ws3.Range("G2").Value = "2016-02-17 03:22:42"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "ffffff9fb0db73-d95f-4e7d-9243-98fa0070a285.md"
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.de-de.xlf"
$ws3.Range("D3").Value = "2016-02-17 03:21:54"
$ws3.Range("E3").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.md"
$ws3.Range("F3").Value = "7ed42eb6-27ad-4a83-baec-771b15d51314.a104ddbcf743322ae326e4e0181e5b36a154e381.de-de.xlf"
$ws3.Range("G3").Value = "2016-02-17 03:22:42"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.34df0f84efa86276470169becf660a08a05a939e.de-de.xlf"
$ws3.Range("D4").Value = "2016-02-17 03:25:54"
$ws3.Range("E4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.md"
$ws3.Range("F4").Value = "fe7f28da-cc97-4606-83ad-7a7b8a3ad11d.34df0f84efa86276470169becf660a08a05a939e.de-de.xlf"
$ws3.Range("G4").Value = "2016-02-17 03:25:00"
$ws3.Range("H4").Value = "Include"
